$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "last updated" timestamp in A1 (04:22 -> 04:52)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 25 de Abril de 2020 a las 04:52"

# ---------------------------------------------------------------------------
# 2) Jamaica gets fresh data and moves (alphabetically) from its old spot
#    (after El Salvador) to right after "Isla de Man" / before "Tanzania".
#    That shifts the Tanzania / Vietnam / El Salvador rows down by one row
#    (rows 123-125 -> 124-126); row 123 becomes Jamaica's updated figures.
# ---------------------------------------------------------------------------
$oldRow123 = $ws.Range("B123:H123").Value2
$oldRow124 = $ws.Range("B124:H124").Value2
$oldRow125 = $ws.Range("B125:H125").Value2

$ws.Range("B126:H126").Value2 = $oldRow125
$ws.Range("B125:H125").Value2 = $oldRow124
$ws.Range("B124:H124").Value2 = $oldRow123

$ws.Range("A123").Value2 = "Jamaica"
$ws.Range("B123").Value2 = 288
$ws.Range("C123").Value2 = 31
$ws.Range("D123").Value2 = 28
$ws.Range("E123").Value2 = 253
$ws.Range("F123").Value2 = 0
$ws.Range("G123").Value2 = 0
$ws.Range("H123").Value2 = 7

$ws.Range("A124").Value2 = "Tanzania"
$ws.Range("A125").Value2 = "Vietnam"
$ws.Range("A126").Value2 = "El Salvador"

# ---------------------------------------------------------------------------
# 3) Haiti (row 159) figures refreshed - no row shift here, straight update.
# ---------------------------------------------------------------------------
$ws.Range("D159").Value2 = 6
$ws.Range("E159").Value2 = 60
$ws.Range("F159").Value2 = 0
$ws.Range("G159").Value2 = 1
$ws.Range("H159").Value2 = 6

Write-Output "edit applied"
